# edit.ps1
# Applies the "final report for group 1" commit to FinalReport.docx using
# Word COM-interop (InsertXML-based run surgery for exact run-level fidelity,
# plus Find for text-only merges).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Locate a Range covering the first occurrence of $text at/after $searchStart.
# Returns a *fresh* Range object (built from plain Long offsets) because
# re-using the Range object that Find.Execute mutated in place causes
# InsertXML to append after the match instead of replacing it.
function Get-RangeForText {
    param($doc, [string]$text, [int]$searchStart = 0)
    $probe = $doc.Range($searchStart, $doc.Content.End)
    $found = $probe.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $text"
    }
    $s = $probe.Start
    $e = $probe.End
    return $doc.Range($s, $e)
}

# Replace the contents of $rng with the literal run/markup XML in $innerXml
# (a fragment of children suitable for sitting directly inside a <w:p>),
# preserving the paragraph's own <w:pPr> and any neighboring paragraphs.
function Set-RunXml {
    param($doc, $rng, [string]$innerXml)
    $pkgXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
              '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
              '<pkg:xmlData>' +
              '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
              '<w:body><w:p>' + $innerXml + '</w:p></w:body></w:document>' +
              '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkgXml)
}

# Insert a brand-new empty paragraph (no runs) with the given <w:pPr> content
# immediately after the paragraph whose text is $afterText.
function Insert-EmptyParagraphAfter {
    param($doc, [string]$afterText, [string]$pprXml)
    $rng = Get-RangeForText $doc $afterText
    $insertAt = $doc.Range($rng.End, $rng.End)
    $pkgXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
              '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
              '<pkg:xmlData>' +
              '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
              '<w:body><w:p><w:pPr>' + $pprXml + '</w:pPr></w:p></w:body></w:document>' +
              '</pkg:xmlData></pkg:part></pkg:package>'
    $insertAt.InsertXML($pkgXml)
}

$stdIndPPr = '<w:ind w:left="1440"/><w:contextualSpacing/>'

# ---------------------------------------------------------------------------
# 1. "W" + "hat you don't..." -> single run (no text change, just a run merge)
# ---------------------------------------------------------------------------
$r = Get-RangeForText $d "What you don’t like about your project and can substantially improve?"
Set-RunXml $d $r '<w:r><w:t>What you don’t like about your project and can substantially improve?</w:t></w:r>'

# ---------------------------------------------------------------------------
# 2. "I Really liked..." -> re-split runs, "Really"->"really", "in the early"->"during the early"
# ---------------------------------------------------------------------------
$r = Get-RangeForText $d "I Really liked the simplicity of it, I wish there would have been more information on how to use it in the early stages of the project."
$xml = '<w:r><w:t xml:space="preserve">I </w:t></w:r>' +
       '<w:r><w:t>r</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">eally liked the simplicity of it, I wish there would have been more information on how to use it </w:t></w:r>' +
       '<w:r><w:t>during</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve"> the early stages of the project.</w:t></w:r>'
Set-RunXml $d $r $xml

# ---------------------------------------------------------------------------
# 3. "It tested well..." -> "It is tested well..." (split into 3 runs)
# ---------------------------------------------------------------------------
$r = Get-RangeForText $d "It tested well for the login and creating new user functions."
$xml = '<w:r><w:t xml:space="preserve">It </w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">is </w:t></w:r>' +
       '<w:r><w:t>tested well for the login and creating new user functions.</w:t></w:r>'
Set-RunXml $d $r $xml

# ---------------------------------------------------------------------------
# 4. "How would you score your confident..." -> merge 2 runs into 1 (no text change)
# ---------------------------------------------------------------------------
$r = Get-RangeForText $d "How would you score your confident in the assessment of the project on a scale of 1 to 5 (1: not confident at all, 5: very confident)? "
Set-RunXml $d $r '<w:r><w:t xml:space="preserve">How would you score your confident in the assessment of the project on a scale of 1 to 5 (1: not confident at all, 5: very confident)? </w:t></w:r>'

# ---------------------------------------------------------------------------
# 5. Insert a new empty paragraph after the confidence-score answer "3"
#    (the one right before "Number and size of commits by each contributor?")
# ---------------------------------------------------------------------------
$r = Get-RangeForText $d "Number and size of commits by each contributor?"
$scoreRng = $d.Range(0, $r.Start)
$scoreRng.Find.Execute("3", $true, $false, $false, $false, $false, $false, 1, $false, "", 0) | Out-Null
$s = $scoreRng.Start
$e = $scoreRng.End
$insertAt = $d.Range($e, $e)
$pkgXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
          '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
          '<pkg:xmlData>' +
          '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body><w:p><w:pPr>' + $stdIndPPr + '</w:pPr></w:p></w:body></w:document>' +
          '</pkg:xmlData></pkg:part></pkg:package>'
$insertAt.InsertXML($pkgXml)

# ---------------------------------------------------------------------------
# 6. James Rodgers commit line: split the "35 commits, 2467+ 1449-" run
# ---------------------------------------------------------------------------
$r = Get-RangeForText $d "35 commits, 2467+ 1449-"
$xml = '<w:r><w:tab/><w:t>35 commits, 2467+</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">+ </w:t></w:r>' +
       '<w:r><w:t>1449-</w:t></w:r>' +
       '<w:r><w:t>-</w:t></w:r>'
Set-RunXml $d $r $xml

# ---------------------------------------------------------------------------
# 7. Yikchun Ng commit line: split the "20 commits 1048+ 536-" run
# ---------------------------------------------------------------------------
$r = Get-RangeForText $d "20 commits 1048+ 536-"
$xml = '<w:r><w:tab/><w:t>20 commits 1048+</w:t></w:r>' +
       '<w:r><w:t>+</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve"> 536-</w:t></w:r>' +
       '<w:r><w:t>-</w:t></w:r>'
Set-RunXml $d $r $xml

# ---------------------------------------------------------------------------
# 8. Napoleon De Mesa commit line: split the "11 commits 146+ 55-" run
# ---------------------------------------------------------------------------
$r = Get-RangeForText $d "11 commits 146+ 55-"
$xml = '<w:r><w:tab/><w:t>11 commits 146+</w:t></w:r>' +
       '<w:r><w:t>+</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve"> 55-</w:t></w:r>' +
       '<w:r><w:t>-</w:t></w:r>'
Set-RunXml $d $r $xml

# ---------------------------------------------------------------------------
# 9. Insert a new empty paragraph after the Napoleon De Mesa commit line
# ---------------------------------------------------------------------------
Insert-EmptyParagraphAfter $d "Napoleon De Mesa:" $stdIndPPr

# ---------------------------------------------------------------------------
# 10. "make  several" (double space, gramStart/gramEnd wrapped) ->
#     "make several" (single space) wrapped by a relocated _GoBack bookmark
# ---------------------------------------------------------------------------
$r = Get-RangeForText $d "make  several"
Set-RunXml $d $r '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>make several</w:t></w:r>'

# ---------------------------------------------------------------------------
# 11. " le" + "ss " -> " less " (merge two runs into one, no text change)
# ---------------------------------------------------------------------------
$r = Get-RangeForText $d " less "
Set-RunXml $d $r '<w:r><w:t xml:space="preserve"> less </w:t></w:r>'

# ---------------------------------------------------------------------------
# 12. Move <w:lastRenderedPageBreak/> from "How would you score the project..."
#     to the start of "How design pattern ..." run
# ---------------------------------------------------------------------------
$r = Get-RangeForText $d "How design pattern "
Set-RunXml $d $r '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">How design pattern </w:t></w:r>'

# ---------------------------------------------------------------------------
# 13. Remove <w:lastRenderedPageBreak/> from "How would you score the project..."
# ---------------------------------------------------------------------------
$r = Get-RangeForText $d "How would you score the project on a scale from 1 to 5?"
Set-RunXml $d $r '<w:r><w:t>How would you score the project on a scale from 1 to 5?</w:t></w:r>'

# ---------------------------------------------------------------------------
# 14. Remove the old _GoBack bookmark (now duplicated after step 10 moved it)
#     from the empty paragraph right after the 2nd "3" score answer.
# ---------------------------------------------------------------------------
$r = Get-RangeForText $d "How many of the issues posted on GitHub project have been resolved?"
$before = $d.Range(0, $r.Start)
$bmFound = $before.Find.Execute("_GoBack_MARKER_NEVER_MATCHES_", $true, $false, $false, $false, $false, $false, 1, $false, "", 0)
# (bookmark text can't be located via Find; locate the empty paragraph by walking Paragraphs instead)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Bookmarks.Count -gt 0) {
        foreach ($bk in $p.Range.Bookmarks) {
            if ($bk.Name -eq "_GoBack") {
                $pStart = $p.Range.Start
                $pEnd = $p.Range.End
                $fresh = $d.Range($pStart, $pEnd - 1)
                if ($fresh.Text -eq "") {
                    Set-RunXml $d $fresh ''
                }
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 15. "In any group activities, members of group " + "experience conflicts..."
#     -> merge into a single run (no text change)
# ---------------------------------------------------------------------------
$r = Get-RangeForText $d "In any group activities, members of group experience conflicts. What was the nature of conflicts you have experienced in this project?"
Set-RunXml $d $r '<w:r><w:t>In any group activities, members of group experience conflicts. What was the nature of conflicts you have experienced in this project?</w:t></w:r>'

# ---------------------------------------------------------------------------
# 16. "Submission: ..." + ".pdf files." -> merge into a single run (no text change)
# ---------------------------------------------------------------------------
$r = Get-RangeForText $d "Submission: There will be two submission pages for this assignment--one for group reports, and one for individual reports. Please submit .pdf files."
Set-RunXml $d $r '<w:r><w:t>Submission: There will be two submission pages for this assignment--one for group reports, and one for individual reports. Please submit .pdf files.</w:t></w:r>'

Write-Host "All edits applied."
